# Generate Report for Archive
#
# The shared status string "Ready for handoff" becomes "In Translation"
# everywhere it is used:
#   - Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3
#   - zh-cn sheet:     column C (Status), rows 2-3
#   - de-de sheet:     column C (Status), rows 2-3
#
# Because the new text is shorter than the old text, the report generator
# also narrows the corresponding status columns (auto-sized to content).

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Update the status text -------------------------------------------------

$wsOverview.Range("E2:F3").Value = $newStatus
$wsZhCn.Range("C2:C3").Value = $newStatus
$wsDeDe.Range("C2:C3").Value = $newStatus

# --- Narrow the status columns to fit the new, shorter text -----------------
# (matches the resized columns the report generator produced)

$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
